$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 411, shifting existing rows 411:459 down to 412:460.
$ws.Rows("411:411").Insert()

# Fill in the new weekly data point in the newly inserted row 411,
# matching the surrounding rows' constant columns (A,B,C,E,F,G,H,I,N,O,Q,R)
# and the new observation's own values for D,J,K,L,M,P.
$ws.Range("A411").Value = 3
$ws.Range("B411").Value = "Femacal de La Calera"
$ws.Range("C411").Value = "Coquimbo"
$ws.Range("D411").Value = 44918
$ws.Range("E411").Value = 5
$ws.Range("F411").Value = 100112009
$ws.Range("G411").Value = "Acelga"
$ws.Range("H411").Value = "Sin especificar"
$ws.Range("I411").Value = "Primera"
$ws.Range("J411").Value = 220
$ws.Range("K411").Value = 3500
$ws.Range("L411").Value = 4000
$ws.Range("M411").Value = 3750
$ws.Range("N411").Value = "`$/docena de atados (6 kilos)"
$ws.Range("O411").Value = "Provincia de Quillota"
$ws.Range("P411").Value = 625
$ws.Range("Q411").Value = 6
$ws.Range("R411").Value = "Hortaliza"
